$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 112, shifting rows 112:119 down to 113:120.
$ws.Rows.Item(112).Insert()

# Populate the newly inserted row 112 with the same "record" metadata as the
# surrounding Sandia / Feria Lagunitas de Puerto Montt rows, and the new
# date / volume / price figures from the weekly update.
$ws.Cells.Item(112, 1).Value = 4
$ws.Cells.Item(112, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(112, 3).Value = "Los Lagos"
$ws.Cells.Item(112, 4).Value = 44476
$ws.Cells.Item(112, 5).Value = 10
$ws.Cells.Item(112, 6).Value = 100112028
$ws.Cells.Item(112, 7).Value = "Sandia"
$ws.Cells.Item(112, 8).Value = "Sin especificar"
$ws.Cells.Item(112, 9).Value = "Primera"
$ws.Cells.Item(112, 10).Value = 450
$ws.Cells.Item(112, 11).Value = 1400
$ws.Cells.Item(112, 12).Value = 1400
$ws.Cells.Item(112, 13).Value = 1400
$ws.Cells.Item(112, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(112, 15).Value = "Per" + [char]0x00FA
$ws.Cells.Item(112, 16).Value = 1400
$ws.Cells.Item(112, 17).Value = 1
$ws.Cells.Item(112, 18).Value = "Hortaliza"

# Match the date-number formatting style used by the rest of column D.
$ws.Cells.Item(112, 4).NumberFormat = $ws.Cells.Item(111, 4).NumberFormat()
